$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.157.86'
$ws.Range("E2").Value = '  -1.48%  '
$ws.Range("D3").Value = '1.573.88'
$ws.Range("E3").Value = '  -0.70%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.41'
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("E6").Value = '  -1.54%  '
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.24'
$ws.Range("E8").Value = '  -0.83%  '
$ws.Range("E9").Value = '  -0.60%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0589'
$ws.Range("E10").Value = '  -0.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0869'
$ws.Range("E11").Value = '  +0.32%  '
$ws.Range("D12").Value = '1.798.82'
$ws.Range("E12").Value = '  -0.64%  '
$ws.Range("D13").Value = '1.582.68'
$ws.Range("E13").Value = '  -0.03%  '
$ws.Range("E14").Value = '  -1.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.518'
$ws.Range("E15").Value = '  -1.49%  '
$ws.Range("D16").Value = '27.191.43'
$ws.Range("E16").Value = '  -1.39%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.24'
$ws.Range("E17").Value = '  -1.46%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.39'
$ws.Range("E18").Value = '  +0.60%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '214.46'
$ws.Range("E19").Value = '  -0.68%  '
$ws.Range("E20").Value = '  -0.95%  '
$ws.Range("E21").Value = '  +0.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.12'
$ws.Range("E22").Value = '  -0.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.48'
$ws.Range("E23").Value = '  -3.30%  '
$ws.Range("E24").Value = '  +0.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.40'
$ws.Range("E25").Value = '  -0.69%  '
$ws.Range("E26").Value = '  -3.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.95'
$ws.Range("E27").Value = '  -0.87%  '
$ws.Range("E28").Value = '  +0.21%  '
$ws.Range("E29").Value = '  -0.92%  '
$ws.Range("E30").Value = '  -3.39%  '
$ws.Range("E31").Value = '  -1.88%  '
$ws.Range("E32").Value = '  -1.33%  '
$ws.Range("D33").Value = '1.398.10'
$ws.Range("E33").Value = '  +1.94%  '
$ws.Range("E34").Value = '  -1.22%  '
$ws.Range("E35").Value = '  +0.61%  '
$ws.Range("E36").Value = '  -0.82%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.942'
$ws.Range("E37").Value = '  -3.37%  '
$ws.Range("E38").Value = '  -2.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.816'
$ws.Range("E39").Value = '  -1.44%  '
$ws.Range("E40").Value = '  -3.62%  '
$ws.Range("E41").Value = '  +0.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.01'
$ws.Range("E42").Value = '  +4.10%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.83'
$ws.Range("E43").Value = '  +1.82%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.43'
$ws.Range("E44").Value = '  +2.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.20'
$ws.Range("E45").Value = '  +1.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '63.81'
$ws.Range("E46").Value = '  -1.16%  '
$ws.Range("D47").Value = '1.710.45'
$ws.Range("E47").Value = '  -0.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.65'
$ws.Range("E48").Value = '  +0.26%  '
$ws.Range("E49").Value = '  -1.25%  '
$ws.Range("E50").Value = '  -0.86%  '
$ws.Range("E51").Value = '  -0.30%  '
